$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 13.36931025828373
$ws.Range("C2").Value = 9.422897799850974
$ws.Range("D2").Value = 6.026864162184811
$ws.Range("E2").Value = 12.73022404250567
$ws.Range("G2").Value = 3.702333941454289
$ws.Range("K2").Value = 9.739424589096744
$ws.Range("L2").Value = 10.02888488763979
$ws.Range("M2").Value = 14.69083017942248
$ws.Range("O2").Value = 31.34136482988521

$ws.Range("B3").Value = 13.20755144935303
$ws.Range("C3").Value = 9.406477123265226
$ws.Range("D3").Value = 5.91527320067409
$ws.Range("E3").Value = 12.75243137404825
$ws.Range("G3").Value = 3.704593855806451
$ws.Range("K3").Value = 9.623372056246186
$ws.Range("L3").Value = 10.03790385804664
$ws.Range("M3").Value = 14.67712137138527
$ws.Range("O3").Value = 31.37618706766923

$ws.Range("B4").Value = 13.11047279115062
$ws.Range("C4").Value = 9.396129532451353
$ws.Range("D4").Value = 5.847468569719559
$ws.Range("E4").Value = 12.76757820744424
$ws.Range("G4").Value = 3.7060552135525
$ws.Range("K4").Value = 9.553658276080297
$ws.Range("L4").Value = 10.04479364810983
$ws.Range("M4").Value = 14.67112766708207
$ws.Range("O4").Value = 31.40308415090342

$ws.Range("B5").Value = 13.071524050185
$ws.Range("C5").Value = 9.391846404576867
$ws.Range("D5").Value = 5.820058271639576
$ws.Range("E5").Value = 12.77413095103176
$ws.Range("G5").Value = 3.706669336677072
$ws.Range("K5").Value = 9.525671219570059
$ws.Range("L5").Value = 10.04794169675047
$ws.Range("M5").Value = 14.66929752719524
$ws.Range("O5").Value = 31.41543054314414

$ws.Range("B6").Value = 13.0650949567113
$ws.Range("C6").Value = 9.391131192635461
$ws.Range("D6").Value = 5.815521351095447
$ws.Range("E6").Value = 12.7752420035887
$ws.Range("G6").Value = 3.706772436999058
$ws.Range("K6").Value = 9.52105044576766
$ws.Range("L6").Value = 10.04848499804331
$ws.Range("M6").Value = 14.66903069676085
$ws.Range("O6").Value = 31.41756429043268

$ws.Range("B7").Value = 13.10994497607706
$ws.Range("C7").Value = 9.396072037152649
$ws.Range("D7").Value = 5.847097956664586
$ws.Range("E7").Value = 12.76766503998942
$ws.Range("G7").Value = 3.706063420408865
$ws.Range("K7").Value = 9.553279080483145
$ws.Range("L7").Value = 10.04483472502188
$ws.Range("M7").Value = 14.67110050221049
$ws.Range("O7").Value = 31.40324505059784

$ws.Range("B8").Value = 13.31309954728628
$ws.Range("C8").Value = 9.417291171357293
$ws.Range("D8").Value = 5.988264483091312
$ws.Range("E8").Value = 12.73756756873
$ws.Range("G8").Value = 3.703097887036277
$ws.Range("K8").Value = 9.69910998448762
$ws.Range("L8").Value = 10.0317142579325
$ws.Range("M8").Value = 14.68560221816337
$ws.Range("O8").Value = 31.35222589646382

$ws.Range("B9").Value = 13.72712544562063
$ws.Range("C9").Value = 9.456785752711767
$ws.Range("D9").Value = 6.268942465225133
$ws.Range("E9").Value = 12.69052915981029
$ws.Range("G9").Value = 3.697865004326369
$ws.Range("K9").Value = 9.995807443501342
$ws.Range("L9").Value = 10.01669417081402
$ws.Range("M9").Value = 14.73312918491212
$ws.Range("O9").Value = 31.29599890814325

$ws.Range("B10").Value = 14.03793161495144
$ws.Range("C10").Value = 9.484498727450518
$ws.Range("D10").Value = 6.475140526334042
$ws.Range("E10").Value = 12.66326117837673
$ws.Range("G10").Value = 3.694371686529809
$ws.Range("K10").Value = 10.21827028257364
$ws.Range("L10").Value = 10.01215943848612
$ws.Range("M10").Value = 14.77947770820453
$ws.Range("O10").Value = 31.28146036359784

$ws.Range("B11").Value = 14.18014400149806
$ws.Range("C11").Value = 9.496820199387743
$ws.Range("D11").Value = 6.568464573774249
$ws.Range("E11").Value = 12.6524361691523
$ws.Range("G11").Value = 3.692857946582139
$ws.Range("K11").Value = 10.32001000742144
$ws.Range("L11").Value = 10.01150078093117
$ws.Range("M11").Value = 14.80299404568719
$ws.Range("O11").Value = 31.28066403968798

$ws.Range("B12").Value = 14.23406370879327
$ws.Range("C12").Value = 9.501444713102551
$ws.Range("D12").Value = 6.603697809457707
$ws.Range("E12").Value = 12.64856381410853
$ws.Range("G12").Value = 3.692295511338366
$ws.Range("K12").Value = 10.35857784060104
$ws.Range("L12").Value = 10.01145254483097
$ws.Range("M12").Value = 14.81224393986339
$ws.Range("O12").Value = 31.28119858948726

$ws.Range("B13").Value = 14.2224489693639
$ws.Range("C13").Value = 9.500450586154422
$ws.Range("D13").Value = 6.596115026592836
$ws.Range("E13").Value = 12.6493877107288
$ws.Range("G13").Value = 3.692416163034784
$ws.Range("K13").Value = 10.35027030842651
$ws.Range("L13").Value = 10.01145399821884
$ws.Range("M13").Value = 14.81023655865999
$ws.Range("O13").Value = 31.28104629013967

$ws.Range("B14").Value = 14.1845789588102
$ws.Range("C14").Value = 9.497201491747738
$ws.Range("D14").Value = 6.571365580051546
$ws.Range("E14").Value = 12.65211304354829
$ws.Range("G14").Value = 3.692811458860354
$ws.Range("K14").Value = 10.32318238692372
$ws.Range("L14").Value = 10.01149278500435
$ws.Range("M14").Value = 14.80374816003338
$ws.Range("O14").Value = 31.28069126467957

$ws.Range("B15").Value = 14.16138968691799
$ws.Range("C15").Value = 9.49520593137045
$ws.Range("D15").Value = 6.55619085557594
$ws.Range("E15").Value = 12.65381192169231
$ws.Range("G15").Value = 3.693054991883598
$ws.Range("K15").Value = 10.3065945481708
$ws.Range("L15").Value = 10.01154271981639
$ws.Range("M15").Value = 14.79981857428531
$ws.Range("O15").Value = 31.28058266637854

$ws.Range("B16").Value = 14.02864986486576
$ws.Range("C16").Value = 9.483687720696201
$ws.Range("D16").Value = 6.469028698710639
$ws.Range("E16").Value = 12.6640003763095
$ws.Range("G16").Value = 3.694472125241399
$ws.Range("K16").Value = 10.21162909029006
$ws.Range("L16").Value = 10.01223067796305
$ws.Range("M16").Value = 14.77798938375607
$ws.Range("O16").Value = 31.28162943824719

$ws.Range("B17").Value = 13.9473925550012
$ws.Range("C17").Value = 9.476548361023443
$ws.Range("D17").Value = 6.415407857914422
$ws.Range("E17").Value = 12.67065498747933
$ws.Range("G17").Value = 3.695360759694863
$ws.Range("K17").Value = 10.15348307556663
$ws.Range("L17").Value = 10.0130118662798
$ws.Range("M17").Value = 14.76521738166602
$ws.Range("O17").Value = 31.28376139015475

$ws.Range("B18").Value = 13.90073620597225
$ws.Range("C18").Value = 9.472415156390968
$ws.Range("D18").Value = 6.38452388401543
$ws.Range("E18").Value = 12.67463121547895
$ws.Range("G18").Value = 3.695878977877923
$ws.Range("K18").Value = 10.12009202624084
$ws.Range("L18").Value = 10.01359340531083
$ws.Range("M18").Value = 14.758100508978
$ws.Range("O18").Value = 31.2855352395617

$ws.Range("B19").Value = 13.88495461895968
$ws.Range("C19").Value = 9.471011125806946
$ws.Range("D19").Value = 6.374060958976813
$ws.Range("E19").Value = 12.6760030401237
$ws.Range("G19").Value = 3.696055658632581
$ws.Range("K19").Value = 10.10879661565036
$ws.Range("L19").Value = 10.01381303479907
$ws.Range("M19").Value = 14.75573037747639
$ws.Range("O19").Value = 31.286229892559

$ws.Range("B20").Value = 13.95603455890749
$ws.Range("C20").Value = 9.477311137675864
$ws.Range("D20").Value = 6.421120587050726
$ws.Range("E20").Value = 12.6699312076431
$ws.Range("G20").Value = 3.695265428708151
$ws.Range("K20").Value = 10.15966759911276
$ws.Range("L20").Value = 10.01291502830312
$ws.Range("M20").Value = 14.76655329208432
$ws.Range("O20").Value = 31.28347776664383

$ws.Range("B21").Value = 14.19570088860556
$ws.Range("C21").Value = 9.498156954700598
$ws.Range("D21").Value = 6.578638278212715
$ws.Range("E21").Value = 12.65130639307426
$ws.Range("G21").Value = 3.692695058625826
$ws.Range("K21").Value = 10.33113793505586
$ws.Range("L21").Value = 10.01147593858945
$ws.Range("M21").Value = 14.80564464434338
$ws.Range("O21").Value = 31.28077285824213

$ws.Range("B22").Value = 14.35270377309241
$ws.Range("C22").Value = 9.511539667513032
$ws.Range("D22").Value = 6.680947809970012
$ws.Range("E22").Value = 12.6404560621649
$ws.Range("G22").Value = 3.691078012279889
$ws.Range("K22").Value = 10.44342750432809
$ws.Range("L22").Value = 10.01170766598032
$ws.Range("M22").Value = 14.83320013854268
$ws.Range("O22").Value = 31.28387797728368

$ws.Range("B23").Value = 14.2688916008052
$ws.Range("C23").Value = 9.504419244887275
$ws.Range("D23").Value = 6.626413685778846
$ws.Range("E23").Value = 12.64612621333262
$ws.Range("G23").Value = 3.691935328781265
$ws.Range("K23").Value = 10.38348783632961
$ws.Range("L23").Value = 10.0114769904638
$ws.Range("M23").Value = 14.81831133513186
$ws.Range("O23").Value = 31.28177509421355

$ws.Range("B24").Value = 13.95212731946155
$ws.Range("C24").Value = 9.476966375704871
$ws.Range("D24").Value = 6.418538035730304
$ws.Range("E24").Value = 12.67025795999886
$ws.Range("G24").Value = 3.695308505008896
$ws.Range("K24").Value = 10.15687145526578
$ws.Range("L24").Value = 10.01295839621275
$ws.Range("M24").Value = 14.76594862288431
$ws.Range("O24").Value = 31.28360428525468

$ws.Range("B25").Value = 13.61374624499343
$ws.Range("C25").Value = 9.446329377577896
$ws.Range("D25").Value = 6.192847548732297
$ws.Range("E25").Value = 12.70197281062477
$ws.Range("G25").Value = 3.69921867338369
$ws.Range("K25").Value = 9.914607037678481
$ws.Range("L25").Value = 10.01961365081095
$ws.Range("M25").Value = 14.7182482647828
$ws.Range("O25").Value = 31.30651101190177

